$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (41-48) appended to the log. All values are written as
# plain text so the exact string representation from the source data is
# preserved (matches the original sheet's convention of storing every
# cell, including numbers, as text).

$rows = @(
    @("2021-04-05","12:15:28","2.2923497267759094","3527.2","3081.6","42.8","42.8","InService"),
    @("2021-04-05","12:22:40","3","3434.0","3087.4","43.8","43.2","InService"),
    @("2021-04-05","12:23:18","2","3416.0","3052.2","43.8","43.2","InService"),
    @("2021-04-05","12:24:08","2.3003","3466.8","3077.8","43.8","43.2","InService"),
    @("2021-04-05","12:24:39","2.3002593312956856","3466.8","3077.8","43.8","43.2","InService"),
    @("2021-04-05","12:25:00","2.3079","3424.0","3045.6","42.4","42.6","InService"),
    @("2021-04-05","12:25:13","2.3079","3424.0","3045.6","42.4","42.6","InService"),
    @("2021-04-05","12:27:28","2.2667","3406.0","3030.4","42.4","42.6","InService")
)

$startRow = 41
$endRow = $startRow + $rows.Length - 1

$fullRange = $ws.Range("A" + $startRow + ":H" + $endRow)
$fullRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
    $ws.Cells.Item($r, 7).Value = $rowData[6]
    $ws.Cells.Item($r, 8).Value = $rowData[7]
}

$fullRange.Style = "Normal"
